$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Calr"
$ws.Range("C2").Value = "Itga3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 93.32574466666667
$ws.Range("H2").Value = 279.977234
$ws.Range("I2").Value = 0.2327963689879921
$ws.Range("J2").Value = 0.2327963689879922
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 12.77892533333333
$ws.Range("N2").Value = 38.336776
$ws.Range("O2").Value = 0.7206984557633653
$ws.Range("P2").Value = 0.7206984557633654
$ws.Range("Q2").Value = 1192.602722773065
$ws.Range("R2").Value = 10733.42450495758
$ws.Range("S2").Value = 0.1677759836369645
$ws.Range("T2").Value = 0.1677759836369646

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Calr"
$ws.Range("C3").Value = "Itga3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 93.32574466666667
$ws.Range("H3").Value = 279.977234
$ws.Range("I3").Value = 0.2327963689879921
$ws.Range("J3").Value = 0.2327963689879922
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.9898276666666668
$ws.Range("N3").Value = 2.969483
$ws.Range("O3").Value = 0.05582372947885773
$ws.Range("P3").Value = 0.05582372947885774
$ws.Range("Q3").Value = 92.37640408333579
$ws.Range("R3").Value = 831.3876367500221
$ws.Range("S3").Value = 0.01299556152604602
$ws.Range("T3").Value = 0.01299556152604602

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Calr"
$ws.Range("C4").Value = "Itga3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 93.32574466666667
$ws.Range("H4").Value = 279.977234
$ws.Range("I4").Value = 0.2327963689879921
$ws.Range("J4").Value = 0.2327963689879922
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02976266666666667
$ws.Range("N4").Value = 0.089288
$ws.Range("O4").Value = 0.001678537697541373
$ws.Range("P4").Value = 0.001678537697541373
$ws.Range("Q4").Value = 2.777623029932445
$ws.Range("R4").Value = 24.998607269392
$ws.Range("S4").Value = 0.0003907574811970962
$ws.Range("T4").Value = 0.0003907574811970963

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Calr"
$ws.Range("C5").Value = "Itga3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 93.32574466666667
$ws.Range("H5").Value = 279.977234
$ws.Range("I5").Value = 0.2327963689879921
$ws.Range("J5").Value = 0.2327963689879922
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.932791
$ws.Range("N5").Value = 11.798373
$ws.Range("O5").Value = 0.2217992770602354
$ws.Range("P5").Value = 0.2217992770602354
$ws.Range("Q5").Value = 367.0306486933646
$ws.Range("R5").Value = 3303.275838240282
$ws.Range("S5").Value = 0.05163406634378446
$ws.Range("T5").Value = 0.05163406634378448

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Calr"
$ws.Range("C6").Value = "Itga3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 164.7897643333334
$ws.Range("H6").Value = 494.369293
$ws.Range("I6").Value = 0.4110597662007076
$ws.Range("J6").Value = 0.4110597662007077
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 12.77892533333333
$ws.Range("N6").Value = 38.336776
$ws.Range("O6").Value = 0.7206984557633653
$ws.Range("P6").Value = 0.7206984557633654
$ws.Range("Q6").Value = 2105.836094113263
$ws.Range("R6").Value = 18952.52484701937
$ws.Range("S6").Value = 0.2962501387273
$ws.Range("T6").Value = 0.2962501387273

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Calr"
$ws.Range("C7").Value = "Itga3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 164.7897643333334
$ws.Range("H7").Value = 494.369293
$ws.Range("I7").Value = 0.4110597662007076
$ws.Range("J7").Value = 0.4110597662007077
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.9898276666666668
$ws.Range("N7").Value = 2.969483
$ws.Range("O7").Value = 0.05582372947885773
$ws.Range("P7").Value = 0.05582372947885774
$ws.Range("Q7").Value = 163.1134679206133
$ws.Range("R7").Value = 1468.021211285519
$ws.Range("S7").Value = 0.02294688918803081
$ws.Range("T7").Value = 0.02294688918803081

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Calr"
$ws.Range("C8").Value = "Itga3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 164.7897643333334
$ws.Range("H8").Value = 494.369293
$ws.Range("I8").Value = 0.4110597662007076
$ws.Range("J8").Value = 0.4110597662007077
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02976266666666667
$ws.Range("N8").Value = 0.089288
$ws.Range("O8").Value = 0.001678537697541373
$ws.Range("P8").Value = 0.001678537697541373
$ws.Range("Q8").Value = 4.904582825931556
$ws.Range("R8").Value = 44.14124543338401
$ws.Range("S8").Value = 0.0006899793135104309
$ws.Range("T8").Value = 0.0006899793135104311

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Calr"
$ws.Range("C9").Value = "Itga3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 164.7897643333334
$ws.Range("H9").Value = 494.369293
$ws.Range("I9").Value = 0.4110597662007076
$ws.Range("J9").Value = 0.4110597662007077
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.932791
$ws.Range("N9").Value = 11.798373
$ws.Range("O9").Value = 0.2217992770602354
$ws.Range("P9").Value = 0.2217992770602354
$ws.Range("Q9").Value = 648.0837020622544
$ws.Range("R9").Value = 5832.753318560289
$ws.Range("S9").Value = 0.09117275897186634
$ws.Range("T9").Value = 0.09117275897186637

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Calr"
$ws.Range("C10").Value = "Itga3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 95.835818
$ws.Range("H10").Value = 287.507454
$ws.Range("I10").Value = 0.2390576204784642
$ws.Range("J10").Value = 0.2390576204784643
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 12.77892533333333
$ws.Range("N10").Value = 38.336776
$ws.Range("O10").Value = 0.7206984557633653
$ws.Range("P10").Value = 0.7206984557633654
$ws.Range("Q10").Value = 1224.678762480923
$ws.Range("R10").Value = 11022.1088623283
$ws.Range("S10").Value = 0.1722884579172938
$ws.Range("T10").Value = 0.1722884579172939

$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Calr"
$ws.Range("C11").Value = "Itga3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 95.835818
$ws.Range("H11").Value = 287.507454
$ws.Range("I11").Value = 0.2390576204784642
$ws.Range("J11").Value = 0.2390576204784643
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.9898276666666668
$ws.Range("N11").Value = 2.969483
$ws.Range("O11").Value = 0.05582372947885773
$ws.Range("P11").Value = 0.05582372947885774
$ws.Range("Q11").Value = 94.86094411403134
$ws.Range("R11").Value = 853.748497026282
$ws.Range("S11").Value = 0.01334508793544923
$ws.Range("T11").Value = 0.01334508793544923

$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Calr"
$ws.Range("C12").Value = "Itga3"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 95.835818
$ws.Range("H12").Value = 287.507454
$ws.Range("I12").Value = 0.2390576204784642
$ws.Range("J12").Value = 0.2390576204784643
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.02976266666666667
$ws.Range("N12").Value = 0.089288
$ws.Range("O12").Value = 0.001678537697541373
$ws.Range("P12").Value = 0.001678537697541373
$ws.Range("Q12").Value = 2.852329505861334
$ws.Range("R12").Value = 25.670965552752
$ws.Range("S12").Value = 0.0004012672278576408
$ws.Range("T12").Value = 0.0004012672278576409

$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Calr"
$ws.Range("C13").Value = "Itga3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 95.835818
$ws.Range("H13").Value = 287.507454
$ws.Range("I13").Value = 0.2390576204784642
$ws.Range("J13").Value = 0.2390576204784643
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 3.932791
$ws.Range("N13").Value = 11.798373
$ws.Range("O13").Value = 0.2217992770602354
$ws.Range("P13").Value = 0.2217992770602354
$ws.Range("Q13").Value = 376.902242508038
$ws.Range("R13").Value = 3392.120182572342
$ws.Range("S13").Value = 0.0530228073978635
$ws.Range("T13").Value = 0.05302280739786351

$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Calr"
$ws.Range("C14").Value = "Itga3"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 46.93870866666666
$ws.Range("H14").Value = 140.816126
$ws.Range("I14").Value = 0.117086244332836
$ws.Range("J14").Value = 0.117086244332836
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 12.77892533333333
$ws.Range("N14").Value = 38.336776
$ws.Range("O14").Value = 0.7206984557633653
$ws.Range("P14").Value = 0.7206984557633654
$ws.Range("Q14").Value = 599.8262532944195
$ws.Range("R14").Value = 5398.436279649776
$ws.Range("S14").Value = 0.08438387548180697
$ws.Range("T14").Value = 0.084383875481807

$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Calr"
$ws.Range("C15").Value = "Itga3"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 46.93870866666666
$ws.Range("H15").Value = 140.816126
$ws.Range("I15").Value = 0.117086244332836
$ws.Range("J15").Value = 0.117086244332836
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.9898276666666668
$ws.Range("N15").Value = 2.969483
$ws.Range("O15").Value = 0.05582372947885773
$ws.Range("P15").Value = 0.05582372947885774
$ws.Range("Q15").Value = 46.46123247587311
$ws.Range("R15").Value = 418.151092282858
$ws.Range("S15").Value = 0.006536190829331674
$ws.Range("T15").Value = 0.006536190829331676

$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Calr"
$ws.Range("C16").Value = "Itga3"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 46.93870866666666
$ws.Range("H16").Value = 140.816126
$ws.Range("I16").Value = 0.117086244332836
$ws.Range("J16").Value = 0.117086244332836
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.02976266666666667
$ws.Range("N16").Value = 0.089288
$ws.Range("O16").Value = 0.001678537697541373
$ws.Range("P16").Value = 0.001678537697541373
$ws.Range("Q16").Value = 1.397021139809778
$ws.Range("R16").Value = 12.573190258288
$ws.Range("S16").Value = 0.0001965336749762051
$ws.Range("T16").Value = 0.0001965336749762052

$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Calr"
$ws.Range("C17").Value = "Itga3"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 46.93870866666666
$ws.Range("H17").Value = 140.816126
$ws.Range("I17").Value = 0.117086244332836
$ws.Range("J17").Value = 0.117086244332836
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 3.932791
$ws.Range("N17").Value = 11.798373
$ws.Range("O17").Value = 0.2217992770602354
$ws.Range("P17").Value = 0.2217992770602354
$ws.Range("Q17").Value = 184.6001309958887
$ws.Range("R17").Value = 1661.401178962998
$ws.Range("S17").Value = 0.0259696443467211
$ws.Range("T17").Value = 0.02596964434672111
